$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Renomear variaveis": the "socio" text that used to live in its own
# column K now lives directly in column J (overwriting the old "Bahia"
# value for that row).
$ws.Range("J5").Value = "46042299000148 - Sócio Pessoa Jurídica Domiciliado no Exterior (Estados Unidos) Representado por Roberto Lazaro dos Santos - Procurador`nRoberto Lazaro dos Santos - Sócio-Administrador"

# The rest of column K (header, blanks, and the K7 socio text) is no
# longer needed, so drop the whole column.
$ws.Columns.Item(11).Delete()
